$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New s_vals data (regenerated to filter save games)
$data = @{
    2 = @{ B = 0.04172184405617529;  C = 0.04103571897497393;  D = 0.7210945179870265;  E = 13.86384647080068;  G = 14.66769855181886 }
    3 = @{ B = 1.445647641019636;    C = 1.626987699542094;    D = 0.1496068669990043;  E = 0.5333859586016987; G = 3.755628166162433 }
    4 = @{ B = 0.1169995834814548;   C = 0.3048912486333797;   D = 0.1496068669990043;  E = 0.5333859586016987; G = 1.104883657715537 }
    5 = @{ B = 0.04172184405617529;  C = 0.00006708468553440206; D = 0.1496068669990043; E = 0.5333859586016987; G = 0.7247817543424127 }
    6 = @{ B = 0.6545652718822623;   C = 0.002658071450198252; D = 0.1496068669990043;  E = 0.5333859586016987; G = 1.340216168933164 }
    7 = @{ B = 3.272327238179451;    C = 1.626987699542094;    D = 3.223369029078222;   E = 0.5333859586016987; G = 8.656069925401464 }
    8 = @{ B = 1.445647641019636;    C = 1.626987699542094;    D = 0.1496068669990043;  E = 0.5333859586016987; G = 3.755628166162433 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
